$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new data record needs to be inserted as row 251, pushing the existing
# rows 251..369 down to 252..370 (one extra row of weekly price data).
# Capture the values that stay the same (unchanged columns) from the row
# that is about to shift down, so they can be reused for the newly
# inserted row 251.
$valA = $ws.Range("A251").Value2
$valB = $ws.Range("B251").Value2
$valC = $ws.Range("C251").Value2
$valE = $ws.Range("E251").Value2
$valF = $ws.Range("F251").Value2
$valG = $ws.Range("G251").Value2
$valH = $ws.Range("H251").Value2
$valI = $ws.Range("I251").Value2
$valN = $ws.Range("N251").Value2
$valQ = $ws.Range("Q251").Value2
$valR = $ws.Range("R251").Value2

# Insert the new row, shifting rows 251:369 down to 252:370.
$ws.Rows("251:251").Insert()

# Populate the newly inserted row 251: unchanged columns copied from the
# captured values above, plus the new weekly record's own data.
$ws.Range("A251").Value2 = $valA
$ws.Range("B251").Value2 = $valB
$ws.Range("C251").Value2 = $valC
$ws.Range("D251").Value2 = 44609
$ws.Range("E251").Value2 = $valE
$ws.Range("F251").Value2 = $valF
$ws.Range("G251").Value2 = $valG
$ws.Range("H251").Value2 = $valH
$ws.Range("I251").Value2 = $valI
$ws.Range("J251").Value2 = 400
$ws.Range("K251").Value2 = 1200
$ws.Range("L251").Value2 = 1200
$ws.Range("M251").Value2 = 1200
$ws.Range("N251").Value2 = $valN
$ws.Range("O251").Value2 = "Provincia de Cautín"
$ws.Range("P251").Value2 = 1200
$ws.Range("Q251").Value2 = $valQ
$ws.Range("R251").Value2 = $valR
